# UltraQueue Benchmarks.xlsx update: add "v.0.7.2" benchmark columns (1 / 32 / 128
# channels) to the Blad1 data sheet, mirroring the existing v.0.7.1 columns
# (H/I, J, K) one step to the right (L, M, N).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Headers (row 2) for the three new series.
# ---------------------------------------------------------------------------
$ws.Range("L2").Value = "v.0.7.2 - 1 Channel"
$ws.Range("M2").Value = "v.0.7.2 - 32 Channels"
$ws.Range("N2").Value = "v.0.7.2 - 128 Channels"

# Copy the header formatting (fill + centered alignment) from the K2 header
# cell onto the three new header cells.
$ws.Range("K2").Copy()
$ws.Range("L2:N2").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2. Data rows 3-21 for the three new series.
# ---------------------------------------------------------------------------
$dataL = @(14.8, 28.8, 53, 98, 205, 403, 739, 1360, 2365, 3213, 4890, 7243, 11870, 12972, 12960, 12496, 12810, 12497, 11669)
$dataM = @(5, 12, 22, 44, 92, 184, 344, 668, 1253, 1916, 3302, 5414, 7957, 10297, 11538, 11795, 12208, 12355, 11417)
$dataN = @(1, 3, 7, 16, 31, 63, 121, 239, 466, 846, 1569, 2831, 4840, 7294, 9051, 10502, 11669, 12073, 11291)

for ($i = 0; $i -lt $dataL.Length; $i++) {
    $row = 3 + $i
    $ws.Cells.Item($row, 12).Value = $dataL[$i]
    $ws.Cells.Item($row, 13).Value = $dataM[$i]
    $ws.Cells.Item($row, 14).Value = $dataN[$i]
}

# Copy the data formatting (centered alignment, no fill) from column K onto
# the new L:N data columns.
$ws.Range("K3:K21").Copy()
$ws.Range("L3:N21").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3. Column widths for the new columns (matching the bestFit widths Excel
#    would compute for the new header text).
# ---------------------------------------------------------------------------
$ws.Columns.Item(12).ColumnWidth = 16.451822916666668
$ws.Columns.Item(13).ColumnWidth = 18.307291666666668
$ws.Columns.Item(14).ColumnWidth = 19.307291666666668

# ---------------------------------------------------------------------------
# 4. Update the sheet view: scroll so column C is left-most and select N22
#    (the cell just below the new data, matching where the author ended up
#    after typing in the new columns).
# ---------------------------------------------------------------------------
$ws.Range("N22").Select()
